$p = $ppt.ActivePresentation
$s = $p.Slides.Add(9, 12)

$tb = $s.Shapes.AddTextbox(1, 0, 0, 720, 54)
$tb.TextFrame.TextRange.Text = "DiSCoVER: top drugs (cerebellar stem cell control)"

$gf = $s.Shapes.AddTable(21, 4, 32.4, 61.2, 651.6, 324.0)
$tbl = $gf.Table
$tbl.Columns.Item(1).Width = 1005840 / 12700
$tbl.Columns.Item(2).Width = 731520 / 12700
$tbl.Columns.Item(3).Width = 1051560 / 12700
$tbl.Columns.Item(4).Width = 5486400 / 12700

$tbl.Rows.Item(1).Height = 195942 / 12700
$tbl.Cell(1,1).Shape.TextFrame.TextRange.Text = "Drug"
$tbl.Cell(1,2).Shape.TextFrame.TextRange.Text = "Score"
$tbl.Cell(1,3).Shape.TextFrame.TextRange.Text = "Evidence"
$tbl.Cell(1,4).Shape.TextFrame.TextRange.Text = "Mechanism of action"

$tbl.Rows.Item(2).Height = 195942 / 12700
$tr = $tbl.Cell(2,1).Shape.TextFrame.TextRange
$tr.Text = "alectinib"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(2,2).Shape.TextFrame.TextRange
$tr.Text = "0.68"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(2,3).Shape.TextFrame.TextRange
$tr.Text = "+.."
$tr.Font.Size = 10.5
$tr = $tbl.Cell(2,4).Shape.TextFrame.TextRange
$tr.Text = "ALK inhibitor, used to treat non-small-cell lung cancer (NSCLC)"
$tr.Font.Size = 10.5

$tbl.Rows.Item(3).Height = 195942 / 12700
$tr = $tbl.Cell(3,1).Shape.TextFrame.TextRange
$tr.Text = "gsk1070916"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(3,2).Shape.TextFrame.TextRange
$tr.Text = "0.64"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(3,3).Shape.TextFrame.TextRange
$tr.Text = "+.."
$tr.Font.Size = 10.5
$tr = $tbl.Cell(3,4).Shape.TextFrame.TextRange
$tr.Text = "Not Clinically Relevant"
$tr.Font.Size = 10.5

$tbl.Rows.Item(4).Height = 195942 / 12700
$tr = $tbl.Cell(4,1).Shape.TextFrame.TextRange
$tr.Text = "sb52334"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(4,2).Shape.TextFrame.TextRange
$tr.Text = "0.62"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(4,3).Shape.TextFrame.TextRange
$tr.Text = "+.."
$tr.Font.Size = 10.5
$tr = $tbl.Cell(4,4).Shape.TextFrame.TextRange
$tr.Text = "Not Clinically Relevant"
$tr.Font.Size = 10.5

$tbl.Rows.Item(5).Height = 195942 / 12700
$tr = $tbl.Cell(5,1).Shape.TextFrame.TextRange
$tr.Text = "ql-xii-61"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(5,2).Shape.TextFrame.TextRange
$tr.Text = "0.59"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(5,3).Shape.TextFrame.TextRange
$tr.Text = "+.."
$tr.Font.Size = 10.5
$tr = $tbl.Cell(5,4).Shape.TextFrame.TextRange
$tr.Text = "Not Clinically Relevant"
$tr.Font.Size = 10.5

$tbl.Rows.Item(6).Height = 195942 / 12700
$tr = $tbl.Cell(6,1).Shape.TextFrame.TextRange
$tr.Text = "tl-2-105"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(6,2).Shape.TextFrame.TextRange
$tr.Text = "0.57"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(6,3).Shape.TextFrame.TextRange
$tr.Text = "+.."
$tr.Font.Size = 10.5
$tr = $tbl.Cell(6,4).Shape.TextFrame.TextRange
$tr.Text = "Not Clinically Relevant"
$tr.Font.Size = 10.5

$tbl.Rows.Item(7).Height = 195942 / 12700
$tr = $tbl.Cell(7,1).Shape.TextFrame.TextRange
$tr.Text = "ql-xi-92"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(7,2).Shape.TextFrame.TextRange
$tr.Text = "0.57"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(7,3).Shape.TextFrame.TextRange
$tr.Text = "+.."
$tr.Font.Size = 10.5
$tr = $tbl.Cell(7,4).Shape.TextFrame.TextRange
$tr.Text = "Not Clinically Relevant"
$tr.Font.Size = 10.5

$tbl.Rows.Item(8).Height = 195942 / 12700
$tr = $tbl.Cell(8,1).Shape.TextFrame.TextRange
$tr.Text = "vx-702"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(8,2).Shape.TextFrame.TextRange
$tr.Text = "0.56"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(8,3).Shape.TextFrame.TextRange
$tr.Text = "+.."
$tr.Font.Size = 10.5
$tr = $tbl.Cell(8,4).Shape.TextFrame.TextRange
$tr.Text = "Not Clinically Relevant"
$tr.Font.Size = 10.5

$tbl.Rows.Item(9).Height = 195942 / 12700
$tr = $tbl.Cell(9,1).Shape.TextFrame.TextRange
$tr.Text = "gsk429286a"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(9,2).Shape.TextFrame.TextRange
$tr.Text = "0.55"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(9,3).Shape.TextFrame.TextRange
$tr.Text = "+.."
$tr.Font.Size = 10.5
$tr = $tbl.Cell(9,4).Shape.TextFrame.TextRange
$tr.Text = "Not Clinically Relevant"
$tr.Font.Size = 10.5

$tbl.Rows.Item(10).Height = 195942 / 12700
$tr = $tbl.Cell(10,1).Shape.TextFrame.TextRange
$tr.Text = "y-39983"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(10,2).Shape.TextFrame.TextRange
$tr.Text = "0.55"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(10,3).Shape.TextFrame.TextRange
$tr.Text = "+.."
$tr.Font.Size = 10.5
$tr = $tbl.Cell(10,4).Shape.TextFrame.TextRange
$tr.Text = "Not Clinically Relevant"
$tr.Font.Size = 10.5

$tbl.Rows.Item(11).Height = 195942 / 12700
$tr = $tbl.Cell(11,1).Shape.TextFrame.TextRange
$tr.Text = "bx-912"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(11,2).Shape.TextFrame.TextRange
$tr.Text = "0.55"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(11,3).Shape.TextFrame.TextRange
$tr.Text = "+.."
$tr.Font.Size = 10.5
$tr = $tbl.Cell(11,4).Shape.TextFrame.TextRange
$tr.Text = "Not Clinically Relevant"
$tr.Font.Size = 10.5

$tbl.Rows.Item(12).Height = 195942 / 12700
$tr = $tbl.Cell(12,1).Shape.TextFrame.TextRange
$tr.Text = "linsitinib"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(12,2).Shape.TextFrame.TextRange
$tr.Text = "0.54"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(12,3).Shape.TextFrame.TextRange
$tr.Text = "++."
$tr.Font.Size = 10.5
$tr = $tbl.Cell(12,4).Shape.TextFrame.TextRange
$tr.Text = "IGF-1R inhibitor"
$tr.Font.Size = 10.5

$tbl.Rows.Item(13).Height = 195942 / 12700
$tr = $tbl.Cell(13,1).Shape.TextFrame.TextRange
$tr.Text = "tubastatin a"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(13,2).Shape.TextFrame.TextRange
$tr.Text = "0.51"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(13,3).Shape.TextFrame.TextRange
$tr.Text = "++."
$tr.Font.Size = 10.5
$tr = $tbl.Cell(13,4).Shape.TextFrame.TextRange
$tr.Text = "Not Clinically Relevant"
$tr.Font.Size = 10.5

$tbl.Rows.Item(14).Height = 195942 / 12700
$tr = $tbl.Cell(14,1).Shape.TextFrame.TextRange
$tr.Text = "gw-2580"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(14,2).Shape.TextFrame.TextRange
$tr.Text = "0.49"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(14,3).Shape.TextFrame.TextRange
$tr.Text = "+.."
$tr.Font.Size = 10.5
$tr = $tbl.Cell(14,4).Shape.TextFrame.TextRange
$tr.Text = "Not Clinically Relevant"
$tr.Font.Size = 10.5

$tbl.Rows.Item(15).Height = 195942 / 12700
$tr = $tbl.Cell(15,1).Shape.TextFrame.TextRange
$tr.Text = "tretinoin"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(15,2).Shape.TextFrame.TextRange
$tr.Text = "0.48"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(15,3).Shape.TextFrame.TextRange
$tr.Text = "++."
$tr.Font.Size = 10.5
$tr = $tbl.Cell(15,4).Shape.TextFrame.TextRange
$tr.Text = "Retinoid"
$tr.Font.Size = 10.5

$tbl.Rows.Item(16).Height = 195942 / 12700
$tr = $tbl.Cell(16,1).Shape.TextFrame.TextRange
$tr.Text = "navitoclax"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(16,2).Shape.TextFrame.TextRange
$tr.Text = "0.48"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(16,3).Shape.TextFrame.TextRange
$tr.Text = "++."
$tr.Font.Size = 10.5
$tr = $tbl.Cell(16,4).Shape.TextFrame.TextRange
$tr.Text = "Bcl-2 family inhibitor: esp Bcl-xL, Bcl-2 and Bcl-w"
$tr.Font.Size = 10.5

$tbl.Rows.Item(17).Height = 195942 / 12700
$tr = $tbl.Cell(17,1).Shape.TextFrame.TextRange
$tr.Text = "kin001-260"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(17,2).Shape.TextFrame.TextRange
$tr.Text = "0.47"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(17,3).Shape.TextFrame.TextRange
$tr.Text = "+.."
$tr.Font.Size = 10.5
$tr = $tbl.Cell(17,4).Shape.TextFrame.TextRange
$tr.Text = "Not Clinically Relevant"
$tr.Font.Size = 10.5

$tbl.Rows.Item(18).Height = 195942 / 12700
$tr = $tbl.Cell(18,1).Shape.TextFrame.TextRange
$tr.Text = "hg-5-88-01"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(18,2).Shape.TextFrame.TextRange
$tr.Text = "0.47"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(18,3).Shape.TextFrame.TextRange
$tr.Text = "+.."
$tr.Font.Size = 10.5
$tr = $tbl.Cell(18,4).Shape.TextFrame.TextRange
$tr.Text = "Not Clinically Relevant"
$tr.Font.Size = 10.5

$tbl.Rows.Item(19).Height = 195942 / 12700
$tr = $tbl.Cell(19,1).Shape.TextFrame.TextRange
$tr.Text = "amuvatinib"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(19,2).Shape.TextFrame.TextRange
$tr.Text = "0.47"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(19,3).Shape.TextFrame.TextRange
$tr.Text = "+.."
$tr.Font.Size = 10.5
$tr = $tbl.Cell(19,4).Shape.TextFrame.TextRange
$tr.Text = "Not Clinically Relevant"
$tr.Font.Size = 10.5

$tbl.Rows.Item(20).Height = 195942 / 12700
$tr = $tbl.Cell(20,1).Shape.TextFrame.TextRange
$tr.Text = "xmd14-99"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(20,2).Shape.TextFrame.TextRange
$tr.Text = "0.47"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(20,3).Shape.TextFrame.TextRange
$tr.Text = "+.."
$tr.Font.Size = 10.5
$tr = $tbl.Cell(20,4).Shape.TextFrame.TextRange
$tr.Text = "Not Clinically Relevant"
$tr.Font.Size = 10.5

$tbl.Rows.Item(21).Height = 195960 / 12700
$tr = $tbl.Cell(21,1).Shape.TextFrame.TextRange
$tr.Text = "xmd13-2"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(21,2).Shape.TextFrame.TextRange
$tr.Text = "0.46"
$tr.Font.Size = 10.5
$tr = $tbl.Cell(21,3).Shape.TextFrame.TextRange
$tr.Text = "+.."
$tr.Font.Size = 10.5
$tr = $tbl.Cell(21,4).Shape.TextFrame.TextRange
$tr.Text = "Not Clinically Relevant"
$tr.Font.Size = 10.5

